$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.035.92'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '1.957.20'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4861'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2946'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06974'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.88%  '
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.34'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').Value = '1.954.01'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07795'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.494'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6991'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '282.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('D17').Value = '31.064.64'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007770'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.62%  '
$ws.Range('D20').Value = '2.209.11'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.509'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.499'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.845'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.19'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.94'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.193'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.392'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.610'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.424'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04934'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7548'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.168'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.734'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02003'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.704'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.523'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '77.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.114'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9035'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '109.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4448'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.131'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.59%  '
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('D48').Value = '1.021.48'
$ws.Range('E48').Value = '  +10.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.346'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1254'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.39%  '
